$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Commit: "Added history starting from 2009-March-31 historical data to all stocks"
# Insert 9 new rows of earlier OHLCV history at row 89 (2019-11-18 .. 2019-11-28),
# which pushes the existing rows 89..161 down to 98..170.
$ws.Range("A89:A97").EntireRow.Insert()

$newRows = @(
    @{ ts = 1574035200; date = "2019-11-18"; o = 1.853; h = 1.86;  l = 1.82;   c = 1.833; vol = 5122497 },
    @{ ts = 1574121600; date = "2019-11-19"; o = 1.847; h = 1.847; l = 1.793;  c = 1.813; vol = 3013648 },
    @{ ts = 1574208000; date = "2019-11-20"; o = 1.813; h = 1.827; l = 1.8;    c = 1.813; vol = 3052798 },
    @{ ts = 1574294400; date = "2019-11-21"; o = 1.82;  h = 1.82;  l = 1.76;   c = 1.76;  vol = 3818998 },
    @{ ts = 1574380800; date = "2019-11-22"; o = 1.76;  h = 1.76;  l = 1.687;  c = 1.74;  vol = 3476998 },
    @{ ts = 1574640000; date = "2019-11-25"; o = 1.74;  h = 1.74;  l = 1.693;  c = 1.713; vol = 1443599 },
    @{ ts = 1574726400; date = "2019-11-26"; o = 1.713; h = 1.827; l = 1.713;  c = 1.8;   vol = 3909448 },
    @{ ts = 1574812800; date = "2019-11-27"; o = 1.82;  h = 1.88;  l = 1.8;    c = 1.847; vol = 5371797 },
    @{ ts = 1574899200; date = "2019-11-28"; o = 1.853; h = 1.873; l = 1.827;  c = 1.84;  vol = 1981499 }
)

$r = 89
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row.ts

    # Columns B (date) and C (id) look numeric/date-like ("2019-11-18", "5292") and
    # would otherwise be auto-converted by COM's usual literal-entry type sniffing;
    # force them to Text, write, then drop back to the default "Normal" style so the
    # cell ends up with the same (unstyled) text content as the rest of the sheet.
    $ws.Cells.Item($r, 2).NumberFormat = "@"
    $ws.Cells.Item($r, 2).Value = $row.date
    $ws.Cells.Item($r, 2).Style = "Normal"

    $ws.Cells.Item($r, 3).NumberFormat = "@"
    $ws.Cells.Item($r, 3).Value = "5292"
    $ws.Cells.Item($r, 3).Style = "Normal"

    $ws.Cells.Item($r, 4).Value = "UWC"
    $ws.Cells.Item($r, 5).Value = $row.o
    $ws.Cells.Item($r, 6).Value = $row.h
    $ws.Cells.Item($r, 7).Value = $row.l
    $ws.Cells.Item($r, 8).Value = $row.c
    $ws.Cells.Item($r, 9).Value = $row.vol

    $r = $r + 1
}
